$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2374.75
$ws.Cells.Item(40, 9).Value = 2583
$ws.Cells.Item(40, 10).Value = 1750
$ws.Cells.Item(40, 11).Value = 2583
$ws.Cells.Item(40, 12).Value = 1750
$ws.Cells.Item(40, 13).Value = -2408
$ws.Cells.Item(40, 14).Value = -2100
$ws.Cells.Item(86, 8).Value = 1517.25
$ws.Cells.Item(86, 9).Value = 1398.3334
$ws.Cells.Item(86, 10).Value = 1874
$ws.Cells.Item(86, 11).Value = 1398.3334
$ws.Cells.Item(86, 12).Value = 1874
$ws.Cells.Item(86, 13).Value = -275.3334
$ws.Cells.Item(86, 14).Value = -4120
$ws.Cells.Item(89, 8).Value = 1517.25
$ws.Cells.Item(89, 9).Value = 1398.3334
$ws.Cells.Item(89, 10).Value = 1874
$ws.Cells.Item(89, 11).Value = 6991.666999999999
$ws.Cells.Item(89, 12).Value = 9370
$ws.Cells.Item(89, 13).Value = -1375.666999999999
$ws.Cells.Item(89, 14).Value = -20602
$ws.Cells.Item(98, 8).Value = 2486.842
$ws.Cells.Item(98, 9).Value = 2569.4443
$ws.Cells.Item(98, 10).Value = 1000
$ws.Cells.Item(98, 11).Value = 2569.4443
$ws.Cells.Item(98, 12).Value = 1000
$ws.Cells.Item(98, 13).Value = -1071.4443
$ws.Cells.Item(98, 14).Value = -3996
$ws.Cells.Item(116, 8).Value = 14558
$ws.Cells.Item(116, 9).Value = 27895
$ws.Cells.Item(116, 11).Value = 27895
$ws.Cells.Item(116, 13).Value = -24453
$ws.Cells.Item(122, 8).Value = 2486.842
$ws.Cells.Item(122, 9).Value = 2569.4443
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 7708.3329
$ws.Cells.Item(122, 12).Value = 3000
$ws.Cells.Item(122, 13).Value = -5258.3329
$ws.Cells.Item(122, 14).Value = -7900
$ws.Cells.Item(132, 8).Value = 1250.2759
$ws.Cells.Item(132, 9).Value = 1157.7037
$ws.Cells.Item(132, 10).Value = 2500
$ws.Cells.Item(132, 11).Value = 3473.1111
$ws.Cells.Item(132, 12).Value = 7500
$ws.Cells.Item(132, 13).Value = -943.1111000000001
$ws.Cells.Item(132, 14).Value = -12560
$ws.Cells.Item(137, 8).Value = 1439.4166
$ws.Cells.Item(137, 9).Value = 1359.3334
$ws.Cells.Item(137, 10).Value = 2000
$ws.Cells.Item(137, 11).Value = 4078.0002
$ws.Cells.Item(137, 12).Value = 6000
$ws.Cells.Item(137, 13).Value = -1528.0002
$ws.Cells.Item(137, 14).Value = -11100
$ws.Cells.Item(138, 8).Value = 2960.6943
$ws.Cells.Item(138, 9).Value = 2621.6296
$ws.Cells.Item(138, 10).Value = 3977.889
$ws.Cells.Item(138, 11).Value = 7864.888800000001
$ws.Cells.Item(138, 12).Value = 11933.667
$ws.Cells.Item(138, 13).Value = -2724.888800000001
$ws.Cells.Item(138, 14).Value = -22213.667

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2114900.8
$ws.Cells.Item(2, 9).Value = 3322873
$ws.Cells.Item(2, 11).Value = 3322873
$ws.Cells.Item(2, 13).Value = -3322760
$ws.Cells.Item(32, 8).Value = 2642.5542
$ws.Cells.Item(32, 9).Value = 1904.6438
$ws.Cells.Item(32, 11).Value = 1904.6438
$ws.Cells.Item(32, 13).Value = -1617.6438
$ws.Cells.Item(61, 8).Value = 3256.75
$ws.Cells.Item(61, 9).Value = 2864.8572
$ws.Cells.Item(61, 11).Value = 2864.8572
$ws.Cells.Item(61, 13).Value = -2652.8572
$ws.Cells.Item(110, 8).Value = 3091.125
$ws.Cells.Item(110, 9).Value = 1840.6
$ws.Cells.Item(110, 10).Value = 5175.3335
$ws.Cells.Item(110, 11).Value = 1840.6
$ws.Cells.Item(110, 12).Value = 5175.3335
$ws.Cells.Item(110, 13).Value = 204.4000000000001
$ws.Cells.Item(110, 14).Value = -9265.333500000001
$ws.Cells.Item(116, 8).Value = 2114900.8
$ws.Cells.Item(116, 9).Value = 3322873
$ws.Cells.Item(116, 11).Value = 3322873
$ws.Cells.Item(116, 13).Value = -3320579
$ws.Cells.Item(132, 8).Value = 1970.3889
$ws.Cells.Item(132, 9).Value = 1128.3334
$ws.Cells.Item(132, 10).Value = 2812.4443
$ws.Cells.Item(132, 11).Value = 3385.0002
$ws.Cells.Item(132, 12).Value = 8437.332900000001
$ws.Cells.Item(132, 13).Value = -855.0001999999999
$ws.Cells.Item(132, 14).Value = -13497.3329
$ws.Cells.Item(136, 8).Value = 3256.75
$ws.Cells.Item(136, 9).Value = 2864.8572
$ws.Cells.Item(136, 11).Value = 8594.571599999999
$ws.Cells.Item(136, 13).Value = -6044.571599999999

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2114900.8
$ws.Cells.Item(3, 9).Value = 3322873
$ws.Cells.Item(3, 11).Value = 3322873
$ws.Cells.Item(3, 13).Value = -3322759
$ws.Cells.Item(134, 8).Value = 25462.777
$ws.Cells.Item(134, 9).Value = 25462.777
$ws.Cells.Item(134, 11).Value = 76388.33099999999
$ws.Cells.Item(134, 13).Value = -73853.33099999999

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(8, 8).Value = 945
$ws.Cells.Item(8, 10).Value = 945
$ws.Cells.Item(8, 12).Value = 945
$ws.Cells.Item(8, 14).Value = -1225
$ws.Cells.Item(22, 8).Value = 799.75
$ws.Cells.Item(22, 9).Value = 399.66666
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 11).Value = 399.66666
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 13).Value = -49.66665999999998
$ws.Cells.Item(22, 14).Value = -2700
$ws.Cells.Item(31, 8).Value = 1942.5385
$ws.Cells.Item(31, 10).Value = 4900
$ws.Cells.Item(31, 12).Value = 4900
$ws.Cells.Item(31, 14).Value = -5490
$ws.Cells.Item(34, 8).Value = 1942.5385
$ws.Cells.Item(34, 10).Value = 4900
$ws.Cells.Item(34, 12).Value = 4900
$ws.Cells.Item(34, 14).Value = -5304
$ws.Cells.Item(58, 8).Value = 4832590
$ws.Cells.Item(58, 9).Value = 6212599.5
$ws.Cells.Item(58, 11).Value = 6212599.5
$ws.Cells.Item(58, 13).Value = -6212396.5
$ws.Cells.Item(99, 8).Value = 2978.5
$ws.Cells.Item(99, 10).Value = 2978.5
$ws.Cells.Item(99, 12).Value = 2978.5
$ws.Cells.Item(99, 14).Value = -5974.5
$ws.Cells.Item(126, 8).Value = 2978.5
$ws.Cells.Item(126, 10).Value = 2978.5
$ws.Cells.Item(126, 12).Value = 8935.5
$ws.Cells.Item(126, 14).Value = -13875.5
$ws.Cells.Item(132, 8).Value = 2229.8518
$ws.Cells.Item(132, 9).Value = 1215.0769
$ws.Cells.Item(132, 10).Value = 3172.1428
$ws.Cells.Item(132, 11).Value = 3645.2307
$ws.Cells.Item(132, 12).Value = 9516.428400000001
$ws.Cells.Item(132, 13).Value = -1115.2307
$ws.Cells.Item(132, 14).Value = -14576.4284
$ws.Cells.Item(134, 8).Value = 1000
$ws.Cells.Item(134, 9).Value = 1000
$ws.Cells.Item(134, 11).Value = 3000
$ws.Cells.Item(134, 13).Value = -465
$ws.Cells.Item(136, 8).Value = 4832590
$ws.Cells.Item(136, 9).Value = 6212599.5
$ws.Cells.Item(136, 11).Value = 18637798.5
$ws.Cells.Item(136, 13).Value = -18635248.5

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 1000000
$ws.Cells.Item(22, 9).Value = 1000000
$ws.Cells.Item(22, 10).Value = 1000000
$ws.Cells.Item(22, 11).Value = 3000000
$ws.Cells.Item(22, 12).Value = 3000000
$ws.Cells.Item(22, 13).Value = -2999831
$ws.Cells.Item(22, 14).Value = -3000338
$ws.Cells.Item(27, 8).Value = 1000000
$ws.Cells.Item(27, 9).Value = 1000000
$ws.Cells.Item(27, 10).Value = 1000000
$ws.Cells.Item(27, 11).Value = 3000000
$ws.Cells.Item(27, 12).Value = 3000000
$ws.Cells.Item(27, 13).Value = -2999898
$ws.Cells.Item(27, 14).Value = -3000204
$ws.Cells.Item(40, 8).Value = 100
$ws.Cells.Item(40, 9).Value = 100
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 400
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -331
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 1700
$ws.Cells.Item(80, 10).Value = 1700
$ws.Cells.Item(80, 12).Value = 5100
$ws.Cells.Item(80, 14).Value = -6972
$ws.Cells.Item(83, 8).Value = 1700
$ws.Cells.Item(83, 10).Value = 1700
$ws.Cells.Item(83, 12).Value = 15300
$ws.Cells.Item(83, 14).Value = -24660
$ws.Cells.Item(122, 8).Value = 869.4167
$ws.Cells.Item(122, 10).Value = 1120.7142
$ws.Cells.Item(122, 12).Value = 10086.4278
$ws.Cells.Item(122, 14).Value = -14986.4278
$ws.Cells.Item(131, 8).Value = 760.96
$ws.Cells.Item(131, 9).Value = 499.25
$ws.Cells.Item(131, 10).Value = 783.7174
$ws.Cells.Item(131, 11).Value = 1497.75
$ws.Cells.Item(131, 12).Value = 2351.1522
$ws.Cells.Item(131, 13).Value = 3542.25
$ws.Cells.Item(131, 14).Value = -12431.1522
$ws.Cells.Item(134, 8).Value = 2116.5454
$ws.Cells.Item(134, 9).Value = 1698.6666
$ws.Cells.Item(134, 11).Value = 5095.9998
$ws.Cells.Item(134, 13).Value = -25.9997999999996
$ws.Cells.Item(136, 8).Value = 3998
$ws.Cells.Item(136, 9).Value = 3998
$ws.Cells.Item(136, 11).Value = 11994
$ws.Cells.Item(136, 13).Value = -6894
$ws.Cells.Item(137, 8).Value = 3113.1304
$ws.Cells.Item(137, 9).Value = 3030
$ws.Cells.Item(137, 10).Value = 3130.6316
$ws.Cells.Item(137, 11).Value = 9090
$ws.Cells.Item(137, 12).Value = 9391.8948
$ws.Cells.Item(137, 13).Value = -3990
$ws.Cells.Item(137, 14).Value = -19591.8948
$ws.Cells.Item(140, 8).Value = 1610.8438
$ws.Cells.Item(140, 10).Value = 2375.7334
$ws.Cells.Item(140, 12).Value = 7127.2002
$ws.Cells.Item(140, 14).Value = -17487.2002

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 6000
$ws.Cells.Item(102, 9).Value = 7499.5
$ws.Cells.Item(102, 10).Value = 3001
$ws.Cells.Item(102, 11).Value = 7499.5
$ws.Cells.Item(102, 12).Value = 3001
$ws.Cells.Item(102, 13).Value = -5877.5
$ws.Cells.Item(102, 14).Value = -6245
$ws.Cells.Item(122, 8).Value = 2787.5715
$ws.Cells.Item(122, 9).Value = 2006.5
$ws.Cells.Item(122, 10).Value = 3100
$ws.Cells.Item(122, 11).Value = 6019.5
$ws.Cells.Item(122, 12).Value = 9300
$ws.Cells.Item(122, 13).Value = -3569.5
$ws.Cells.Item(122, 14).Value = -14200
$ws.Cells.Item(126, 8).Value = 1716209.9
$ws.Cells.Item(126, 9).Value = 2060078.8
$ws.Cells.Item(126, 11).Value = 6180236.4
$ws.Cells.Item(126, 13).Value = -6177766.4
$ws.Cells.Item(132, 8).Value = 5497904.5
$ws.Cells.Item(132, 9).Value = 5497904.5
$ws.Cells.Item(132, 11).Value = 16493713.5
$ws.Cells.Item(132, 13).Value = -16491183.5

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2065.5386
$ws.Cells.Item(46, 9).Value = 1100
$ws.Cells.Item(46, 10).Value = 2241.0908
$ws.Cells.Item(46, 11).Value = 1100
$ws.Cells.Item(46, 12).Value = 2241.0908
$ws.Cells.Item(46, 13).Value = -912
$ws.Cells.Item(46, 14).Value = -2617.0908

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 544.2414
$ws.Cells.Item(107, 9).Value = 431.2
$ws.Cells.Item(107, 11).Value = 1293.6
$ws.Cells.Item(107, 13).Value = 626.4000000000001
$ws.Cells.Item(132, 8).Value = 1567.7894
$ws.Cells.Item(132, 9).Value = 1102.9286
$ws.Cells.Item(132, 10).Value = 2869.4
$ws.Cells.Item(132, 11).Value = 3308.7858
$ws.Cells.Item(132, 12).Value = 8608.200000000001
$ws.Cells.Item(132, 13).Value = -778.7857999999997
$ws.Cells.Item(132, 14).Value = -13668.2
$ws.Cells.Item(135, 8).Value = 75757.82000000001
$ws.Cells.Item(135, 10).Value = 75757.82000000001
$ws.Cells.Item(135, 12).Value = 75757.82000000001
$ws.Cells.Item(135, 14).Value = -85897.82000000001
$ws.Cells.Item(137, 8).Value = 95000
$ws.Cells.Item(137, 10).Value = 95000
$ws.Cells.Item(137, 12).Value = 95000
$ws.Cells.Item(137, 14).Value = -105200
